# Applies the commit's change to StructureDefinition-Specialite.xlsx:
#   - "Metadata" sheet, cell B5 ("Title" row) gets the same value as
#     B4 ("Name" row): "Specialite".
#   - The regenerated "Date" metadata value (B8) is refreshed to the
#     new generation timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Title (row 5) was blank; the regeneration filled it in with the same
# value as Name (row 4) -> "Specialite".
$ws.Range("B5").Value = $ws.Range("B4").Value2

# Date (row 8) gets bumped to the new generation timestamp.
$ws.Range("B8").Value = "2025-07-17T14:35:50+00:00"
